# Auto-generated edit script: update Leve profit calculation cells
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the
# scheduled price-data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1718.5
$ws.Range("I19").Value = 1549.6666
$ws.Range("J19").Value = 2225
$ws.Range("K19").Value = 1549.6666
$ws.Range("L19").Value = 2225
$ws.Range("M19").Value = -1374.6666
$ws.Range("N19").Value = -2575
$ws.Range("H40").Value = 1520.1
$ws.Range("I40").Value = 1400.2
$ws.Range("J40").Value = 1640
$ws.Range("K40").Value = 1400.2
$ws.Range("L40").Value = 1640
$ws.Range("M40").Value = -1225.2
$ws.Range("N40").Value = -1990
$ws.Range("H43").Value = 8920.625
$ws.Range("I43").Value = 7998.1
$ws.Range("J43").Value = 10458.167
$ws.Range("K43").Value = 7998.1
$ws.Range("L43").Value = 10458.167
$ws.Range("M43").Value = -7929.1
$ws.Range("N43").Value = -10596.167
$ws.Range("H55").Value = 392.4
$ws.Range("J55").Value = 467
$ws.Range("L55").Value = 467
$ws.Range("N55").Value = -895
$ws.Range("H70").Value = 1037
$ws.Range("I70").Value = 937
$ws.Range("J70").Value = 1062
$ws.Range("K70").Value = 2811
$ws.Range("L70").Value = 3186
$ws.Range("M70").Value = -2541
$ws.Range("N70").Value = -3726
$ws.Range("H73").Value = 1037
$ws.Range("I73").Value = 937
$ws.Range("J73").Value = 1062
$ws.Range("K73").Value = 2811
$ws.Range("L73").Value = 3186
$ws.Range("M73").Value = -1875
$ws.Range("N73").Value = -5058
$ws.Range("H76").Value = 4075.5
$ws.Range("I76").Value = 4075.5
$ws.Range("K76").Value = 4075.5
$ws.Range("M76").Value = -3760.5
$ws.Range("H79").Value = 4075.5
$ws.Range("I79").Value = 4075.5
$ws.Range("K79").Value = 4075.5
$ws.Range("M79").Value = -2983.5
$ws.Range("H86").Value = 7787.5
$ws.Range("I86").Value = 6766.6665
$ws.Range("K86").Value = 6766.6665
$ws.Range("M86").Value = -5643.6665
$ws.Range("H89").Value = 7787.5
$ws.Range("I89").Value = 6766.6665
$ws.Range("K89").Value = 33833.3325
$ws.Range("M89").Value = -28217.3325
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null
$ws.Range("H106").Value = 4966.6665
$ws.Range("I106").Value = 4966.6665
$ws.Range("K106").Value = 4966.6665
$ws.Range("M106").Value = -4335.6665
$ws.Range("H107").Value = 829.86664
$ws.Range("I107").Value = 829.86664
$ws.Range("K107").Value = 829.86664
$ws.Range("M107").Value = 1090.13336
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -1558
$ws.Range("H132").Value = 1996.05
$ws.Range("I132").Value = 1738.625
$ws.Range("J132").Value = 3025.75
$ws.Range("K132").Value = 5215.875
$ws.Range("L132").Value = 9077.25
$ws.Range("M132").Value = -2685.875
$ws.Range("N132").Value = -14137.25
$ws.Range("H137").Value = 2319.8
$ws.Range("I137").Value = 1900
$ws.Range("K137").Value = 5700
$ws.Range("M137").Value = -3150

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2493.7273
$ws.Range("I88").Value = 1093
$ws.Range("K88").Value = 1093
$ws.Range("M88").Value = -687
$ws.Range("H91").Value = 2493.7273
$ws.Range("I91").Value = 1093
$ws.Range("K91").Value = 1093
$ws.Range("M91").Value = 311
$ws.Range("H132").Value = 370.81818
$ws.Range("I132").Value = 307.9
$ws.Range("K132").Value = 923.6999999999999
$ws.Range("M132").Value = 1606.3

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2372.647
$ws.Range("I94").Value = 1277.9166
$ws.Range("K94").Value = 1277.9166
$ws.Range("M94").Value = -826.9166
$ws.Range("H105").Value = 3940.5
$ws.Range("I105").Value = 3578
$ws.Range("K105").Value = 3578
$ws.Range("M105").Value = -1831

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1679.4445
$ws.Range("I31").Value = 1622
$ws.Range("J31").Value = 1966.6666
$ws.Range("K31").Value = 1622
$ws.Range("L31").Value = 1966.6666
$ws.Range("M31").Value = -1327
$ws.Range("N31").Value = -2556.6666
$ws.Range("H34").Value = 1679.4445
$ws.Range("I34").Value = 1622
$ws.Range("J34").Value = 1966.6666
$ws.Range("K34").Value = 1622
$ws.Range("L34").Value = 1966.6666
$ws.Range("M34").Value = -1420
$ws.Range("N34").Value = -2370.6666
$ws.Range("H58").Value = 1356.579
$ws.Range("I58").Value = 985.06665
$ws.Range("K58").Value = 985.06665
$ws.Range("M58").Value = -782.06665
$ws.Range("H62").Value = 2837.8
$ws.Range("I62").Value = 2797.25
$ws.Range("K62").Value = 2797.25
$ws.Range("M62").Value = -2173.25
$ws.Range("H65").Value = 2837.8
$ws.Range("I65").Value = 2797.25
$ws.Range("K65").Value = 13986.25
$ws.Range("M65").Value = -10866.25
$ws.Range("H99").Value = 7477.353
$ws.Range("I99").Value = 7721.5713
$ws.Range("J99").Value = 6337.6665
$ws.Range("K99").Value = 7721.5713
$ws.Range("L99").Value = 6337.6665
$ws.Range("M99").Value = -6223.5713
$ws.Range("N99").Value = -9333.666499999999
$ws.Range("H105").Value = 7925.2666
$ws.Range("I105").Value = 10298.091
$ws.Range("J105").Value = 1400
$ws.Range("K105").Value = 10298.091
$ws.Range("L105").Value = 1400
$ws.Range("M105").Value = -8551.091
$ws.Range("N105").Value = -4894
$ws.Range("H126").Value = 7477.353
$ws.Range("I126").Value = 7721.5713
$ws.Range("J126").Value = 6337.6665
$ws.Range("K126").Value = 23164.7139
$ws.Range("L126").Value = 19012.9995
$ws.Range("M126").Value = -20694.7139
$ws.Range("N126").Value = -23952.9995
$ws.Range("H136").Value = 1356.579
$ws.Range("I136").Value = 985.06665
$ws.Range("K136").Value = 2955.19995
$ws.Range("M136").Value = -405.1999500000002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 509.2857
$ws.Range("I122").Value = 134.5
$ws.Range("K122").Value = 1210.5
$ws.Range("M122").Value = 1239.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 13921.75
$ws.Range("I107").Value = 1995.5
$ws.Range("K107").Value = 1995.5
$ws.Range("M107").Value = -75.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1584.5
$ws.Range("J22").Value = 2605
$ws.Range("L22").Value = 2605
$ws.Range("N22").Value = -3195
$ws.Range("H27").Value = 1584.5
$ws.Range("J27").Value = 2605
$ws.Range("L27").Value = 2605
$ws.Range("N27").Value = -2819
$ws.Range("H93").Value = 1519.2
$ws.Range("I93").Value = 1099.7142
$ws.Range("J93").Value = 2498
$ws.Range("K93").Value = 1099.7142
$ws.Range("L93").Value = 2498
$ws.Range("M93").Value = 148.2858000000001
$ws.Range("N93").Value = -4994
$ws.Range("H122").Value = 1646.5
$ws.Range("I122").Value = 1565.6
$ws.Range("J122").Value = 1781.3334
$ws.Range("K122").Value = 4696.799999999999
$ws.Range("L122").Value = 5344.0002
$ws.Range("M122").Value = -2246.799999999999
$ws.Range("N122").Value = -10244.0002
$ws.Range("H132").Value = 1753.7778
$ws.Range("I132").Value = 1660.5
$ws.Range("K132").Value = 4981.5
$ws.Range("M132").Value = -2451.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = $null
$ws.Range("H122").Value = 1535.8182
$ws.Range("I122").Value = 1535.8182
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4607.4546
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2157.4546
$ws.Range("N122").Value = $null
$ws.Range("H132").Value = 1858.4375
$ws.Range("I132").Value = 1925.6666
$ws.Range("J132").Value = 850
$ws.Range("K132").Value = 5776.9998
$ws.Range("L132").Value = 2550
$ws.Range("M132").Value = -3246.9998
$ws.Range("N132").Value = -7610
